$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns O and P, matching style of N1 (bold header style)
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column O values (rows 2-7)
$ws.Range("O2").Value = -1.282010457652405
$ws.Range("O3").Value = -0.4626403529148906
$ws.Range("O4").Value = 0.05704889946834561
$ws.Range("O5").Value = 0.448172929842699
$ws.Range("O6").Value = -0.06680551348815365
$ws.Range("O7").Value = -0.006914363506040602

# Column P values (rows 2-7)
$ws.Range("P2").Value = -1.127490096490725
$ws.Range("P3").Value = -0.4172087221114515
$ws.Range("P4").Value = 0.03204889941959649
$ws.Range("P5").Value = 0.4270283932206927
$ws.Range("P6").Value = -0.06043743938363743
$ws.Range("P7").Value = -0.006555558134232935
